$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellD = $ws.Cells.Item(2, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "21.205.11"
$cellD.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +3.92%  "

$cellD = $ws.Cells.Item(3, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "1.540.73"
$cellD.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +5.38%  "

$cellD = $ws.Cells.Item(4, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "1.003"
$cellD.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.56%  "

$cellD = $ws.Cells.Item(5, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.9600"
$cellD.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +1.02%  "

$cellD = $ws.Cells.Item(6, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "281.65"
$cellD.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +2.63%  "

$cellD = $ws.Cells.Item(7, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.3621"
$cellD.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -0.82%  "

$cellD = $ws.Cells.Item(8, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.3184"
$cellD.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +3.51%  "

$cellD = $ws.Cells.Item(9, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "40.89"
$cellD.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +3.20%  "

$cellD = $ws.Cells.Item(10, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "1.100"
$cellD.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +6.30%  "

$cellD = $ws.Cells.Item(11, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.06800"
$cellD.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +3.59%  "

$cellD = $ws.Cells.Item(12, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.9971"
$cellD.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -0.26%  "

$cellD = $ws.Cells.Item(13, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "5.656"
$cellD.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +4.50%  "

$cellD = $ws.Cells.Item(14, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "18.75"
$cellD.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +4.45%  "

$cellD = $ws.Cells.Item(15, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "6.337"
$cellD.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +3.62%  "

$cellD = $ws.Cells.Item(16, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.00001047"
$cellD.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +2.22%  "

$cellD = $ws.Cells.Item(17, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.9607"
$cellD.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -0.60%  "

$cellD = $ws.Cells.Item(18, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "1.535.06"
$cellD.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +5.05%  "

$cellD = $ws.Cells.Item(19, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.06032"
$cellD.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +4.47%  "

$cellD = $ws.Cells.Item(20, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "71.98"
$cellD.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +3.34%  "

$cellD = $ws.Cells.Item(21, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "5.673"
$cellD.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +4.58%  "

$cellD = $ws.Cells.Item(22, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "15.06"
$cellD.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +4.29%  "

$ws.Cells.Item(23, 5).Value = "  +4.41%  "

$cellD = $ws.Cells.Item(24, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "2.312"
$cellD.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +3.29%  "

$cellD = $ws.Cells.Item(25, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "21.267.99"
$cellD.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +4.07%  "

$cellD = $ws.Cells.Item(26, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "147.62"
$cellD.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +4.34%  "

$cellD = $ws.Cells.Item(27, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "2.200"
$cellD.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +5.75%  "

$cellD = $ws.Cells.Item(28, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "17.71"
$cellD.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +3.50%  "

$cellD = $ws.Cells.Item(29, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "1.700.89"
$cellD.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +5.41%  "

$cellD = $ws.Cells.Item(30, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "117.75"
$cellD.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +5.13%  "

$cellD = $ws.Cells.Item(31, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "4.035"
$cellD.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +5.04%  "

$cellD = $ws.Cells.Item(32, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.8491"
$cellD.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +7.72%  "

$cellD = $ws.Cells.Item(33, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "5.180"
$cellD.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +6.45%  "

$cellD = $ws.Cells.Item(34, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.08027"
$cellD.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +2.81%  "

$cellD = $ws.Cells.Item(35, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "1.500"
$cellD.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -0.66%  "

$cellD = $ws.Cells.Item(36, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "1.216"
$cellD.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +7.21%  "

$cellD = $ws.Cells.Item(37, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "4.955"
$cellD.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +6.38%  "

$cellD = $ws.Cells.Item(38, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.05862"
$cellD.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +3.03%  "

$cellD = $ws.Cells.Item(39, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.02095"
$cellD.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +3.51%  "

$ws.Cells.Item(40, 5).Value = "  +3.77%  "

$cellD = $ws.Cells.Item(41, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "7.691"
$cellD.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +3.23%  "

$cellD = $ws.Cells.Item(42, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.1912"
$cellD.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +3.01%  "

$cellD = $ws.Cells.Item(43, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.9604"
$cellD.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +0.52%  "

$cellD = $ws.Cells.Item(44, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.5455"
$cellD.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +3.81%  "

$cellD = $ws.Cells.Item(45, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "12.42"
$cellD.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +4.76%  "

$cellD = $ws.Cells.Item(46, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "3.564"
$cellD.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +2.28%  "

$cellD = $ws.Cells.Item(47, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.5437"
$cellD.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +6.04%  "

$cellD = $ws.Cells.Item(48, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "121.48"
$cellD.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +3.86%  "

$cellD = $ws.Cells.Item(49, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "1.869"
$cellD.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +7.04%  "

$cellD = $ws.Cells.Item(50, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.06628"
$cellD.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +3.27%  "

$cellD = $ws.Cells.Item(51, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "70.26"
$cellD.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +6.21%  "
